$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.687.85'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +2.76%  '
$ws.Range('E2').Style = "Normal"

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.432.30'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +3.05%  '
$ws.Range('E3').Style = "Normal"

# Row 4
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E4').Style = "Normal"

# Row 5
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('E5').Style = "Normal"

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '130.45'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +3.66%  '
$ws.Range('E6').Style = "Normal"

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.596'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.96%  '
$ws.Range('E7').Style = "Normal"

# Row 8
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E8').Style = "Normal"

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.694'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +5.53%  '
$ws.Range('E9').Style = "Normal"

# Row 10
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +16.03%  '
$ws.Range('E10').Style = "Normal"

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.99'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +2.42%  '
$ws.Range('E11').Style = "Normal"

# Row 12
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('E12').Style = "Normal"

# Row 13
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +1.80%  '
$ws.Range('E13').Style = "Normal"

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '19.83'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +2.89%  '
$ws.Range('E14').Style = "Normal"

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.446.34'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +4.16%  '
$ws.Range('E15').Style = "Normal"

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '62.698.03'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +2.98%  '
$ws.Range('E16').Style = "Normal"

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '11.50'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +2.56%  '
$ws.Range('E17').Style = "Normal"

# Row 18
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('E18').Style = "Normal"

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0000158'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +23.76%  '
$ws.Range('E19').Style = "Normal"

# Row 20
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('E20').Style = "Normal"

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '84.60'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +5.49%  '
$ws.Range('E21').Style = "Normal"

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '314.98'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +5.15%  '
$ws.Range('E22').Style = "Normal"

# Row 23
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E23').Style = "Normal"

# Row 24
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +2.26%  '
$ws.Range('E24').Style = "Normal"

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.76'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +2.28%  '
$ws.Range('E25').Style = "Normal"

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '29.74'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +2.82%  '
$ws.Range('E26').Style = "Normal"

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.17'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('E27').Style = "Normal"

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.84'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +6.15%  '
$ws.Range('E28').Style = "Normal"

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.76'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +8.87%  '
$ws.Range('E29').Style = "Normal"

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '44.52'
$ws.Range('D30').Style = "Normal"

# Row 31
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.70%  '
$ws.Range('E31').Style = "Normal"

# Row 32
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('E32').Style = "Normal"

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '11.40'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E33').Style = "Normal"

# Row 34
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E34').Style = "Normal"

# Row 35
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +1.96%  '
$ws.Range('E35').Style = "Normal"

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '51.88'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('E36').Style = "Normal"

# Row 37
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('E37').Style = "Normal"

# Row 38
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.96'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +2.18%  '
$ws.Range('E38').Style = "Normal"

# Row 39
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.323'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +15.53%  '
$ws.Range('E39').Style = "Normal"

# Row 40
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('E40').Style = "Normal"

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '142.79'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +5.57%  '
$ws.Range('E41').Style = "Normal"

# Row 42
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +2.35%  '
$ws.Range('E42').Style = "Normal"

# Row 43
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.76%  '
$ws.Range('E43').Style = "Normal"

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.90'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +1.57%  '
$ws.Range('E44').Style = "Normal"

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.92'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +1.19%  '
$ws.Range('E45').Style = "Normal"

# Row 46
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('E46').Style = "Normal"

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '21.25'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('E47').Style = "Normal"

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.106.06'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('E48').Style = "Normal"

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.98'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +6.31%  '
$ws.Range('E49').Style = "Normal"

# Row 50
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.45%  '
$ws.Range('E50').Style = "Normal"

# Row 51
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +31.78%  '
$ws.Range('E51').Style = "Normal"
